$d = $word.ActiveDocument

# Locate the anchor paragraph "Ragionare su questa cosa " and expand the
# range to cover the whole paragraph (incl. paragraph mark), then take the
# position right after it as the insertion point for the new paragraphs.
$anchor = $d.Content
$found = $anchor.Find.Execute("Ragionare su questa cosa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Expand(4)
$insPoint = $anchor.End
$insRange = $d.Range($insPoint, $insPoint)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$fragment = '<w:p ' + $ns + '><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p>' +
  '<w:p ' + $ns + '><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>MG - 29/05/2023</w:t></w:r></w:p>' +
  '<w:p ' + $ns + '><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>FATTO:</w:t></w:r></w:p>' +
  '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">Ho inserito nel </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>db</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> i 72 pz che avevano l’infezione sul sito ma non sul batterio.</w:t></w:r></w:p>' +
  '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Ho aggiunto al batterio anche questi nuovi valori</w:t></w:r></w:p>' +
  '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">Ho calcolato </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>num_infezioni</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>batterio_pos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> considerando questi nuovi dati</w:t></w:r></w:p>' +
  '<w:p ' + $ns + '><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>'

$insRange.InsertXML($fragment)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
